$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2027972027972028
$ws.Range("C2").Value = 0.5244755244755245
$ws.Range("J2").Value = 0.02097902097902098
$ws.Range("P2").Value = 0.1328671328671329
$ws.Range("S2").Value = 0.1188811188811189
# Row 3
$ws.Range("B3").Value = 0.006578947368421052
$ws.Range("C3").Value = 0.02631578947368421
$ws.Range("J3").Value = 0.01973684210526316
$ws.Range("P3").Value = 0.6776315789473685
$ws.Range("S3").Value = 0.2697368421052632
# Row 4
$ws.Range("J4").Value = 0.06666666666666667
$ws.Range("P4").Value = 0.6888888888888889
$ws.Range("S4").Value = 0.2444444444444444
# Row 6
$ws.Range("B6").Value = 0.03864734299516908
$ws.Range("D6").Value = 0.02415458937198068
$ws.Range("F6").Value = 0.02898550724637681
$ws.Range("J6").Value = 0.2270531400966184
$ws.Range("O6").Value = 0.01932367149758454
$ws.Range("Q6").Value = 0.2077294685990338
$ws.Range("R6").Value = 0.06280193236714976
$ws.Range("S6").Value = 0.391304347826087
# Row 7
$ws.Range("B7").Value = 0.08536585365853659
$ws.Range("D7").Value = 0.03658536585365853
$ws.Range("F7").Value = 0.06097560975609756
$ws.Range("J7").Value = 0.1158536585365854
$ws.Range("O7").Value = 0.02439024390243903
$ws.Range("Q7").Value = 0.1524390243902439
$ws.Range("R7").Value = 0.09146341463414634
$ws.Range("S7").Value = 0.4329268292682927
# Row 8
$ws.Range("B8").Value = 0.08641975308641975
$ws.Range("D8").Value = 0.00823045267489712
$ws.Range("F8").Value = 0.08436213991769548
$ws.Range("J8").Value = 0.09259259259259259
$ws.Range("O8").Value = 0.01440329218106996
$ws.Range("Q8").Value = 0.1790123456790123
$ws.Range("R8").Value = 0.05967078189300411
$ws.Range("S8").Value = 0.4753086419753086
# Row 9
$ws.Range("B9").Value = 0.1241379310344828
$ws.Range("D9").Value = 0.03448275862068965
$ws.Range("F9").Value = 0.05517241379310345
$ws.Range("J9").Value = 0.1103448275862069
$ws.Range("O9").Value = 0.006896551724137931
$ws.Range("Q9").Value = 0.1724137931034483
$ws.Range("R9").Value = 0.06896551724137931
$ws.Range("S9").Value = 0.4275862068965517
# Row 10
$ws.Range("B10").Value = 0.1151020408163265
$ws.Range("D10").Value = 0.02204081632653061
$ws.Range("E10").Value = 0.0008163265306122449
$ws.Range("F10").Value = 0.0783673469387755
$ws.Range("J10").Value = 0.1012244897959184
$ws.Range("O10").Value = 0.01224489795918367
$ws.Range("Q10").Value = 0.2269387755102041
$ws.Range("R10").Value = 0.06693877551020408
$ws.Range("S10").Value = 0.3763265306122449
# Row 11
$ws.Range("G11").Value = 0.1654135338345865
$ws.Range("J11").Value = 0.09398496240601503
$ws.Range("K11").Value = 0.2142857142857143
$ws.Range("L11").Value = 0.5112781954887218
$ws.Range("S11").Value = 0.01503759398496241
# Row 12
$ws.Range("G12").Value = 0.7971014492753623
$ws.Range("J12").Value = 0.1666666666666667
$ws.Range("K12").Value = 0.01449275362318841
$ws.Range("L12").Value = 0.01449275362318841
$ws.Range("S12").Value = 0.007246376811594203
# Row 13
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.3125
$ws.Range("S13").Value = 0.0625
# Row 15
$ws.Range("F15").Value = 0.01470588235294118
$ws.Range("H15").Value = 0.2352941176470588
$ws.Range("I15").Value = 0.0392156862745098
$ws.Range("J15").Value = 0.3823529411764706
$ws.Range("K15").Value = 0.06862745098039216
$ws.Range("M15").Value = 0.004901960784313725
$ws.Range("O15").Value = 0.05392156862745098
$ws.Range("S15").Value = 0.2009803921568628
# Row 16
$ws.Range("F16").Value = 0.02380952380952381
$ws.Range("H16").Value = 0.244047619047619
$ws.Range("I16").Value = 0.05357142857142857
$ws.Range("J16").Value = 0.4404761904761905
$ws.Range("K16").Value = 0.06547619047619048
$ws.Range("O16").Value = 0.06547619047619048
$ws.Range("S16").Value = 0.1071428571428571
# Row 17
$ws.Range("F17").Value = 0.00881057268722467
$ws.Range("H17").Value = 0.1894273127753304
$ws.Range("I17").Value = 0.08370044052863436
$ws.Range("J17").Value = 0.4625550660792951
$ws.Range("K17").Value = 0.0881057268722467
$ws.Range("M17").Value = 0.01101321585903084
$ws.Range("O17").Value = 0.05947136563876652
$ws.Range("S17").Value = 0.09691629955947137
# Row 18
$ws.Range("F18").Value = 0.006802721088435374
$ws.Range("H18").Value = 0.217687074829932
$ws.Range("I18").Value = 0.08843537414965986
$ws.Range("J18").Value = 0.4557823129251701
$ws.Range("K18").Value = 0.08843537414965986
$ws.Range("M18").Value = 0.006802721088435374
$ws.Range("O18").Value = 0.06122448979591837
$ws.Range("S18").Value = 0.07482993197278912
# Row 19
$ws.Range("F19").Value = 0.01076158940397351
$ws.Range("H19").Value = 0.2301324503311258
$ws.Range("I19").Value = 0.06456953642384106
$ws.Range("J19").Value = 0.4048013245033112
$ws.Range("K19").Value = 0.1043046357615894
$ws.Range("M19").Value = 0.02152317880794702
$ws.Range("O19").Value = 0.07119205298013245
$ws.Range("S19").Value = 0.09271523178807947
